$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C11").Value = 182048
$ws.Range("C19").Select()
